# The deck's primary theme (ppt/theme/theme1.xml, currently named "Integral")
# is being swapped for the stock "Office Theme" palette. fontScheme and
# fmtScheme are already byte-identical between the two themes, so the only
# observable difference is the 10 non-shared colours in <a:clrScheme>
# (dk1/lt1 -- black/white -- are common to both and left untouched).
#
# VBA/COM's ColorScheme.Colors(i).RGB uses the classic 0x00BBGGRR packing
# (same as the VBA RGB() function), so convert each target hex swatch
# (RRGGBB, as it appears in the OOXML <a:srgbClr val="..."/>) into that
# integer form before assigning it.

function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target "Office Theme" colour scheme, in clrScheme document order.
$officeTheme = [ordered]@{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

foreach ($index in $officeTheme.Keys) {
    $colorScheme.Colors($index).RGB = ConvertTo-BgrInt $officeTheme[$index]
}
